$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 25) with the next date and error count,
# reusing the date style already used by the rest of column A.
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A25").Value = 45992
$ws.Range("B25").Value = 2

# Match the selection state recorded in the saved workbook.
$ws.Range("A25:B25").Select()
